# The workbook tracks upcoming "漫展" (comic con) events. The event that
# used to occupy row 2 ("丽水·YA●怀旧only", 2024.02.14) has passed, so it is
# dropped from the sheet; the event that used to be row 3
# ("丽水·LPJ 现实X次元动漫展（取消）", 2024.02.18) moves up to become the new
# row 2 (keeping row 2's original serial-number id in column A), and the old
# row 3 is removed entirely, shrinking the sheet back down to 2 rows.
#
# This happened identically on both the "展览" (sheet 1) and "全部类型"
# (sheet 4) tabs, which mirror each other's data.

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # Overwrite row 2 (B:I) with what used to be row 3's values. Column A
    # (the running index, value 1) is left untouched.
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024.02.18"
    $ws.Range("C2").Value = "丽水·LPJ 现实X次元动漫展（取消）"
    $ws.Range("D2").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E2").Value = "2024.02.18 10:00-02.18 17:00"
    $ws.Range("F2").Value = 314
    $ws.Range("G2").NumberFormat = "@"
    $ws.Range("G2").Value = "不可售"
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=79437"
    $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202312/ee5hLUN61702276208812.jpeg"

    # The old row 3 (now duplicated into row 2) is removed entirely, which
    # shrinks the sheet's used range/dimension down to A1:I2.
    $ws.Rows.Item(3).Delete()
}
